$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the milestone / date cells that were showing raw python dict/JSON text ---
# Row 2: milestone 1 / date
$ws.Range("A2").Value = "Milestone 1"

# B2 needs to hold the literal text "2022-01-15" (not an Excel date serial).
# Enter it as a formula returning text, then freeze it to a static value via
# Copy / PasteSpecial(values) so no NumberFormat/style gets minted and no
# date auto-conversion happens.
$ws.Range("B2").Formula = '="2022-01-15"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Row 3: milestone 3, no completion date anymore
$ws.Range("A3").Value = "Milestone 3"
$ws.Range("B3").ClearContents()

# Row 4 (new): milestone 5, no completion date
$ws.Range("A4").Value = "Milestone 5"

# --- Resize the table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B4"))
